$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.391.99"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.821.15"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'315.58"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.5247"
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("D8").Value = "'0.3854"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("D9").Value = "'0.08054"
$ws.Range("E9").Value = "  +5.31%  "
$ws.Range("D10").Value = "'41.85"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "'1.113"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'6.397"
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "'20.86"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "'7.426"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "1.825.34"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'94.50"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").Value = "'0.00001103"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'17.63"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'6.026"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").Value = "28.434.51"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'11.38"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").Value = "'2.247"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'159.06"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("D27").Value = "'20.87"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "2.025.87"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'2.413"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "'124.45"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("D32").Value = "'1.079"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").Value = "'5.671"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Value = "'3.680"
$ws.Range("D35").Value = "'0.07326"
$ws.Range("E35").Value = "  +4.14%  "
$ws.Range("D36").Value = "'12.17"
$ws.Range("E36").Value = "  +8.54%  "
$ws.Range("D37").Value = "'0.2202"
$ws.Range("D38").Value = "'0.02344"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "'5.127"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'8.752"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "'0.6308"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "'1.180"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").Value = "'1.383"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.49"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6124"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").Value = "'3.794"
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("D47").Value = "'127.24"
$ws.Range("E47").Value = "  +1.67%  "
$ws.Range("D48").Value = "'1.221"
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").Value = "'1.969"
$ws.Range("D50").Value = "'0.06897"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Value = "'73.98"
$ws.Range("E51").Value = "  +0.00%  "
